# UseCase_GestioneOrdini.docx - "Correzione di alcune exit condition"
#
# This script re-creates, via Word COM-interop Range.InsertXML, the set of
# run-splitting / proofErr (spell-check marker) edits described by the
# target diff. Range.InsertXML replaces the *entire enclosing paragraph* in
# this runtime, so every replacement below carries the paragraph's original
# w:pPr and paragraph-mark attributes forward unchanged, and only the run
# content inside is restructured the way the diff shows.

$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Find-NthOccurrence($searchText, $n) {
    $r = $d.Content.Duplicate
    $r.Start = 0
    for ($i = 0; $i -lt $n; $i++) {
        $found = $r.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $found) {
            Write-Host "NOT FOUND occurrence $n of: $searchText"
            return $null
        }
        if ($i -lt $n - 1) {
            $r.Collapse(0)
            $r.End = $d.Content.End
        }
    }
    return $r
}

function Replace-Paragraph($searchText, $occurrence, $newParaXml) {
    $r = Find-NthOccurrence $searchText $occurrence
    if ($null -eq $r) { return }
    $r.InsertXML($pkgOpen + $newParaXml + $pkgClose)
}

# 1. "Vers." -> "Vers" + proofErr(spellStart/spellEnd) + "."
Replace-Paragraph "Vers." 1 '<w:p w14:paraId="3C9859EB" w14:textId="77777777" w:rsidR="00AA34DF" w:rsidRPr="00506DC0" w:rsidRDefault="00AA34DF" w:rsidP="00AB0E6B"><w:pPr><w:rPr><w:bCs/><w:i/><w:iCs/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>Vers</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>.</w:t></w:r></w:p>'

# 2. "Entry Condition" -> "Entry " + proofErr("Condition")
Replace-Paragraph "Entry Condition" 1 '<w:p w14:paraId="25F49BBD" w14:textId="77777777" w:rsidR="009048BC" w:rsidRDefault="00506DC0" w:rsidP="00AB0E6B"><w:r w:rsidRPr="00506DC0"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Entry </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Condition</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'

# 3. First "Exit condition" (under "On success") -> "Exit " + proofErr("condition")
Replace-Paragraph "Exit condition" 1 '<w:p w14:paraId="4AAE6F6F" w14:textId="77777777" w:rsidR="009048BC" w:rsidRDefault="009048BC" w:rsidP="00AB0E6B"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="009048BC"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Exit </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>condition</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'

# 4. Second "Exit condition" (under "On failure") -> "Exit " + proofErr("condition")
Replace-Paragraph "Exit condition" 2 '<w:p w14:paraId="0821136D" w14:textId="77777777" w:rsidR="00E12504" w:rsidRDefault="00E12504" w:rsidP="006F6523"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="009048BC"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Exit </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>condition</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'

# 5. "On failure" paragraph -> "On " + proofErr("failure")
Replace-Paragraph "On failure" 1 '<w:p w14:paraId="77F2E50E" w14:textId="77777777" w:rsidR="00E12504" w:rsidRDefault="00E12504" w:rsidP="00E12504"><w:r><w:t xml:space="preserve">                       On </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>failure</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'

# 6. " non riesce a comunicare col server e non riesce a vedere la pagina." ->
#    " non riesce a " + "visualizzare" + " la pagina" + " " + italic "che contiene l'elenco degli ordini effettuati dai clienti."
Replace-Paragraph "non riesce a comunicare" 1 '<w:p w14:paraId="3C4A4D96" w14:textId="01B58FE7" w:rsidR="00E12504" w:rsidRPr="00506DC0" w:rsidRDefault="00E41F43" w:rsidP="006F6523"><w:r><w:rPr><w:bCs/></w:rPr><w:t>Il consulente</w:t></w:r><w:r w:rsidR="00D95CCD"><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve"> non riesce a </w:t></w:r><w:r><w:rPr><w:bCs/></w:rPr><w:t>visualizzare</w:t></w:r><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve"> la pagina</w:t></w:r><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>che contiene l’elenco degli ordini effettuati dai clienti.</w:t></w:r></w:p>'

# 7. "/User Priority" -> "/User " + proofErr("Priority")
Replace-Paragraph "/User Priority" 1 '<w:p w14:paraId="58549609" w14:textId="77777777" w:rsidR="009048BC" w:rsidRDefault="009048BC" w:rsidP="00AB0E6B"><w:r w:rsidRPr="009048BC"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Rilevanza</w:t></w:r><w:r w:rsidR="00E12504"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">/User </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Priority</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'

# 8. "Generalization of" -> proofErr("Generalization") + " of"
Replace-Paragraph "Generalization of" 1 '<w:p w14:paraId="6660BFF3" w14:textId="77777777" w:rsidR="00F14455" w:rsidRPr="009455E7" w:rsidRDefault="00F14455" w:rsidP="006F6523"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:b/><w:bCs/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="009455E7"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:b/><w:bCs/></w:rPr><w:t>Generalization</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> of</w:t></w:r></w:p>'

# 9. "Special Requirements" -> "Special " + proofErr("Requirements")
Replace-Paragraph "Special Requirements" 1 '<w:p w14:paraId="45637F5C" w14:textId="006D5838" w:rsidR="00CB1A10" w:rsidRDefault="00CB1A10" w:rsidP="000A5648"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Special </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Requirements</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'

Write-Host "Done."
